# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ordering: for each period (ascending: 2003..2012, 2101..2106),
# list worker ISMAEL GUILLERMO APONTE MARIMON (CC 73073941) then
# worker ENRIQUE CARLOS HURTADO COHEN (CC 9157124).
# Salario Basico (col G) updated from 877804 to 877805 for all rows.

$data = @(
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2003",33125,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2003",37472,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2004",33125,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2004",37472,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2005",33125,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2005",37472,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2006",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2006",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2007",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2007",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2008",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2008",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2009",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2009",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2010",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2010",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2011",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2011",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2012",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2012",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2101",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2101",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2102",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2102",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2103",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2103",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2104",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2104",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2105",35112,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2105",35112,877805),
    @("CC","73073941","ISMAEL GUILLERMO APONTE MARIMON","2106",24578,877805),
    @("CC","9157124","ENRIQUE CARLOS HURTADO COHEN","2106",24578,877805)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rec[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]   # G - Salario Basico
}
